# coursework.docx edit ("added spec to course work"):
#   1. Prepend a new paragraph with the module / student header line
#      ("Information -Society and Security" ... "40161070" ... "Adam Blance").
#   2. Turn the former first paragraph's "3." into "1. Robot" (keeping the
#      _GoBack bookmark right after it) and append " designs", " ", " " as
#      three more runs in that same paragraph.
#   3. Insert a fresh paragraph containing just "3." right after that, before
#      the existing "4." paragraph.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new paragraph ("Information -Society and Security"
# plus the tab-separated module id / author line) in front of the very first
# paragraph in the document (the one that currently just holds "3." and the
# _GoBack bookmark). Inserting a full <w:p> fragment right at the very start
# of the body (offset 0) creates a clean, standalone new paragraph without
# disturbing what follows.
# ---------------------------------------------------------------------------
$docStart = $d.Range(0, 0)

$headerXml = @"
<w:p xmlns:w="$wNs">
  <w:r><w:t>Information -Society and Security</w:t></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/><w:t>40161070</w:t></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/><w:t xml:space="preserve">Adam </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Blance</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/></w:r>
</w:p>
"@

$docStart.InsertXML($headerXml)

# ---------------------------------------------------------------------------
# Step 2: the paragraph that used to be first ("3." + _GoBack bookmark) is now
# the second paragraph. Turn its "3." into "1. Robot" (the bookmark, which
# sits right after the run in the XML, is left untouched / simply follows the
# replaced text).
# ---------------------------------------------------------------------------
$listPara = $d.Paragraphs(2).Range
$listPara.Find.Execute("3.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Robot", 2)

# ---------------------------------------------------------------------------
# Step 3: append " designs", " " and " " (three separate runs) at the end of
# that same paragraph, after the bookmark.
#
# A <w:p>-wrapped InsertXML fragment always carries its own paragraph break,
# so inserting it at the shared boundary between paragraph 2 and paragraph 3
# merges the new runs onto the *front* of paragraph 3 instead of the end of
# paragraph 2. We exploit that (it keeps the three runs distinct instead of
# being coalesced into one run the way Range.InsertAfter would), then split
# paragraph 3 again right before "4." and finally delete the paragraph mark
# that still separates paragraph 2 from the new runs, re-merging them into a
# single paragraph together with the "1. Robot" text and the bookmark.
# ---------------------------------------------------------------------------
$listPara = $d.Paragraphs(2).Range
$boundary = $d.Range($listPara.End, $listPara.End)
$tailXml = @"
<w:p xmlns:w="$wNs">
  <w:r><w:t xml:space="preserve"> designs</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
"@
$boundary.InsertXML($tailXml)

# Paragraph 3 now reads " designs   4." - split it again right before "4." so
# the tail runs end up alone in their own paragraph.
$mergedPara = $d.Paragraphs(3).Range
$fourFinder = $mergedPara.Duplicate
$fourFinder.Find.Execute("4.")
$splitPoint = $d.Range($fourFinder.Start, $fourFinder.Start)
$splitPoint.InsertParagraphBefore()

# Re-merge paragraph 2 ("1. Robot" + bookmark) and paragraph 3 (the tail
# runs) into one paragraph by deleting the paragraph mark between them.
$listPara = $d.Paragraphs(2).Range
$mergeMark = $d.Range($listPara.End - 1, $listPara.End)
$mergeMark.Delete()

# ---------------------------------------------------------------------------
# Step 4: insert a new paragraph containing just "3." right after that
# paragraph (i.e. before the paragraph "4."). Just like in step 3, inserting
# the <w:p> fragment at the paragraph boundary merges its runs onto the front
# of the "4." paragraph, so split that paragraph again right before "4." to
# leave "3." on its own.
# ---------------------------------------------------------------------------
$listPara = $d.Paragraphs(2).Range
$afterListPara = $d.Range($listPara.End, $listPara.End)
$threeXml = @"
<w:p xmlns:w="$wNs">
  <w:r><w:t>3.</w:t></w:r>
</w:p>
"@
$afterListPara.InsertXML($threeXml)

$mergedPara2 = $d.Paragraphs(3).Range
$fourFinder2 = $mergedPara2.Duplicate
$fourFinder2.Find.Execute("4.")
$splitPoint2 = $d.Range($fourFinder2.Start, $fourFinder2.Start)
$splitPoint2.InsertParagraphBefore()

Write-Host "edit complete"
